# edit.ps1 - apply "Only 14 genomes left!" commit
# Fills in newly-completed metabolic/functional annotation columns for
# genomes in rows 159-174 (J, L, N, P, Q, S, T), and corrects six rows where
# Motility data had been typed into the wrong column (R -> S).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 159-174: newly filled annotation cells ---
# Row 159
$ws.Range("J159").Value = 'chitobiose, MurNAc, glucose, glycolate, 1,3-B-glucan, glucoside, cellobiose, starch/glycogen, maltose, trehalose'
$ws.Range("N159").Value = 'nitrate reduction, nitrogen fixation, ammonia_assimilation, hydroxylamine reductase'
$ws.Range("P159").Value = 'sulfate_red_dis, trithionate'
$ws.Range("Q159").Value = 'branched amino, amino acid/amide, cobalt/nickel, L-amino, iron, LPS export, lipoprotein release, microcin C, multiple sugar, molybdate, phosphate, phospholipid/cholesterol, phosphonate, tungstate'
$ws.Range("S159").Value = 'chemotaxis, flagellum'
$ws.Range("T159").Value = 'Oxidative phosphorylation'

# Row 160
$ws.Range("J160").Value = 'three carotenoid genes, chitobiose, fructose, glucose, galactose, sucrose, glucoside, cellulose, cellobiose, starch/glycogen, trehalsoe, maltose, dextrin, isomaltose'
$ws.Range("N160").Value = 'partial denitrification, ammonia_assimilation'
$ws.Range("Q160").Value = 'Cu-processing, heme, LPS export, lipoprotein release, molybdate, phosphate, phospholipid/cholesterol'
$ws.Range("S160").Value = 'one chemotaxis protein'
$ws.Range("T160").Value = 'Oxidative phosphorylation'

# Row 161
$ws.Range("J161").Value = 'chitobiose, fructose, starch/glycogen, trehalose'
$ws.Range("N161").Value = 'ammonia_assimilation'
$ws.Range("Q161").Value = 'peptidase, arginine, biotin, LPS export, lipoprotein release, macrolide, phospholipid/cholesterol, putative hydroxymethylpyrimidine, spermidine/putrescine, type IV secretion'
$ws.Range("T161").Value = 'Oxidative phosphorylation'

# Row 162
$ws.Range("J162").Value = 'chitin, chitobiose, MurNAc, fructose, glucose, sucrose, stachyose, raffinose, glycolate, glycerate, starch/glycogen, maltose'
$ws.Range("N162").Value = 'nitrate_red_diss, partial denitrification, ammonia_assimilation'
$ws.Range("P162").Value = 'sufate_red_ass, thiosulfate, tetrathionate'
$ws.Range("Q162").Value = 'amino acid/amide, branched amino, cobalt/nickel, iron, LPS export, lipoprotein release, maltose/maltodextrin, molybdate, multiple sugar, phosphate, phospholipid/cholesterol, sodium, tungstate, type VI secretion'
$ws.Range("S162").Value = 'chemotaxis (purine?), flagellum'
$ws.Range("T162").Value = 'Oxidative phosphorylation'

# Row 163
$ws.Range("J163").Value = 'chitobiose, glucose, galactose, sorbose, sorbitol, rhamnose, rhamnulose, raffinose/stachyose/manninotriose, melibiose, tartrate, glycolate, trehalose, galacturonate,  cellulose, starch/glycogen, maltose'
$ws.Range("N163").Value = 'nitrate reduction, ammonia_assimilation'
$ws.Range("P163").Value = 'sulfate_red_ass, trithionate'
$ws.Range("Q163").Value = 'branched amino, MacB-like domain, NAG, amino acid/amide, carbohydrate, cobalt/nickel, extracellular, iron, iron (III), lactose/arabinose, LPS export, lipoprotein release, molybdate, monosaccharide, multiple sugar, oligopeptide, phosphate, phospholipid/cholesterol, raffinose/stachyose/melibiose, sodium, xylose, type IV secretion'
$ws.Range("T163").Value = 'Oxidative phosphorylation'

# Row 164
$ws.Range("J164").Value = 'partial Wood-Ljungdahl, acetate, formaldehyde, formate'
$ws.Range("L164").Value = 'chitobiose, glucose, tartrate, glycolate, glycerate, cellulose, galacturonate, starch/glycogen, dextrin, trehalose, maltose'
$ws.Range("N164").Value = 'nitrate_red_dis, partial denitrification, nitrogen fixation, ammonia_assimilation'
$ws.Range("P164").Value = 'sulfate_red_ass, thiosulfate'
$ws.Range("Q164").Value = 'branched amino, polysaccharide/polyol phosphate export, cobalt, permease, zinc, amino acid/amide, cobalamin, cobalt/nickel, LPS export, lipoprotein release, macrolide, molybdate, molybdenum, phosphate, phospholipid/cholesterol, sulfate, tungstate, zinc, type II secretion, type VI secretion'
$ws.Range("S164").Value = 'chemotaxis (purine?), flagellum'
$ws.Range("T164").Value = 'Oxidative phosphorylation'

# Row 165
$ws.Range("J165").Value = 'two carotenoid genes, galactose, fructan, sucrose, cellobiose, cellulose synthesis, maltose'
$ws.Range("N165").Value = 'partial denitrification, ammonia_assimilation'
$ws.Range("Q165").Value = 'Cu-processing, permease, LPS export, lipoprotein release, phosphate, phospholipid/cholesterol'
$ws.Range("S165").Value = 'one chemotaxis protein'
$ws.Range("T165").Value = 'Oxidative phosphorylation'

# Row 166
$ws.Range("J166").Value = 'chitobiose, fructose, galactose, fructan, sorbitol, rhamnulose, raffinose, stachyose, manninotriose, meliobiose, sucrose, galactan, glucoside, cellulose, cellobiose, starch/glycogen, pectin, maltose'
$ws.Range("N166").Value = 'nitrite reducatse, nitroalkane, ammonia_assimilation, hydroxylamine reductase'
$ws.Range("P166").Value = 'sulfate_red_ass, tetrathionate'
$ws.Range("Q166").Value = 'cobalt, iron, LPS export, lipoprotein, oligopeptide/dipeptide, phospholipid/cholesterol, sulfate'
$ws.Range("S166").Value = 'chemotaxis (purine?)'
$ws.Range("T166").Value = 'Oxidative phosphorylation'

# Row 167
$ws.Range("J167").Value = 'chitobiose, fructose, glucose, fructan, rhamnose, rhamnulose, galactose, galactan, lactose, raffinose, stachyose, manninotriose, melibiose,,pectin, cellobiose, glucoside, starch/glycogen, maltose'
$ws.Range("N167").Value = 'nitrogen_fixation, ammonia_assimilation, hydroxylamine reductase'
$ws.Range("P167").Value = 'alkanesulfonate, methanesulfonate, sulfite reductaase'
$ws.Range("Q167").Value = 'permease, LPS export, lipoprotein release, molybdate, phosphate, phospholipid, ribose, sulfonate, ribose'
$ws.Range("S167").Value = 'one chemotaxis protein'
$ws.Range("T167").Value = 'Oxidative phosphorylation, luciferase-like monooxygenase'

# Row 168
$ws.Range("Q168").Value = 'LPS transport'
$ws.Range("T168").Value = 'Oxidative phosphorylation'

# Row 169
$ws.Range("J169").Value = 'partial Wood-Ljungdahl, acetate, formate'
$ws.Range("L169").Value = 'chitobiose, glucose, starch/glycogen'
$ws.Range("N169").Value = 'nitrite reductase, ammonia_assimilation'
$ws.Range("P169").Value = 'sulfate_red_ass, thiosulfate'
$ws.Range("Q169").Value = 'MacB-like domain, permease, heme, LPS export, lipoprotein, molybdate, phosphate, phospholipid/cholesterol'
$ws.Range("S169").Value = 'a few chemotaxis protein'
$ws.Range("T169").Value = 'Oxidative phosphorylation'

# Row 170
$ws.Range("J170").Value = 'one carotenoid gene, formate, partial Wood-Ljungdahl, acetate'
$ws.Range("L170").Value = 'chitin, chitobiose, fructose, glucose, galactose, fructan, fucose , fuculose, rhamnose, rhamnulose, lactose, galactan, raffinose, stachyose, manninotriose, melibiose, sucrose, pectin, glucoside, cellobiose, maltose, trehalose'
$ws.Range("N170").Value = 'nitrite reductase, nitric oxide reductase, ammonia_assimilation'
$ws.Range("P170").Value = 'thiosulfate'
$ws.Range("Q170").Value = 'NitT/TauT, permease, heme, iron, LPS export, lipoprotein release, macrolide, molybdate, peptide/nickel, phosphate, phospholipid/cholesterol, phosphonate, tungstate/molybdate, zinc'
$ws.Range("S170").Value = 'chemotaxis (purine?), flagellum'
$ws.Range("T170").Value = 'Oxidative phosphorylation'

# Row 171
$ws.Range("J171").Value = 'partial Wood-Ljungdahl, acetate, formate'
$ws.Range("L171").Value = 'chitin, chitobiose, glucosaminide, fructose, galactan, lactose, sucrose, glycerate, cellobiose, glucoside, starch/glycogen, maltose'
$ws.Range("N171").Value = 'nitrite reductase, nitric oxide reductase, ammonia_assimilation'
$ws.Range("P171").Value = 'sulfate_red_ass, alkanesulfonate'
$ws.Range("Q171").Value = 'NitT/TauT, biotin, branched amino, carbohydrate, heme, LPS export, lipoprotein release, molybdate, phosphate, phospholipid/cholesterol, zinc'
$ws.Range("S171").Value = 'one chemotaxis protein'
$ws.Range("T171").Value = 'Oxidative phosphorylation'

# Row 172
$ws.Range("J172").Value = 'chitobiose, MurNAc, glucoside, cellobiose'
$ws.Range("N172").Value = 'ammonia_assimilation'
$ws.Range("Q172").Value = 'biotin, heme, LPS export, lipoprotein release, phospholipid/cholesterol, spermidine/putrescine, zinc, type IV secretion'
$ws.Range("T172").Value = 'Oxidative phosphorylation'

# Row 173
$ws.Range("J173").Value = 'chitobiose, MurNAc'
$ws.Range("N173").Value = 'ammonia_assimilation'
$ws.Range("Q173").Value = 'LPS export, lipoprotein release, oligopeptide, phospholipid/cholesterol'
$ws.Range("S173").Value = 'a few chemotaxis proteins, flagellum'
$ws.Range("T173").Value = 'Oxidative phosphorylation'

# Row 174
$ws.Range("J174").Value = 'Some carotenoid genes, partial Wood-Ljungdahl, acetate, formate'
$ws.Range("L174").Value = 'chitobiose, fructose, galactose, galactonate, tartrate, glycerate, glycolate'
$ws.Range("N174").Value = 'nitrate reductase, nitrite reductase, ammonia_assimilation, formamide, nitroalkane, nitrile'
$ws.Range("P174").Value = 'sulfate oxidation, taurine, alkanesulfonate, methanesulfonate, thiosulfate'
$ws.Range("Q174").Value = 'HCOMODA decarboxylase, amino acid/amide, glutamate, branched amino, L-amino, glycine betaine/choline, heme, iron, iron(III), LPS export, lipoprotein release, molybdate, multiple sugar, nitrate/nitrite, phosphate, phospholipid/cholesterol, polar amino, sulfonate, tungstate, urea'
$ws.Range("S174").Value = 'a few chemotaxis and flagellum proteins'
$ws.Range("T174").Value = 'Oxidative phosphorylation'

# --- Rows 190,193,194,197,201,203: Motility values were mistakenly entered
# in column R (Annual_Trend); move them to column S (Motility). ---
$ws.Range("S190").Value = 'one chemotaxis protein'
$ws.Range("R190").ClearContents()
$ws.Range("S193").Value = 'chemotaxis (purine?), flagellum'
$ws.Range("R193").ClearContents()
$ws.Range("S194").Value = 'one chemotaxis protein'
$ws.Range("R194").ClearContents()
$ws.Range("S197").Value = 'two flagellum proteins'
$ws.Range("R197").ClearContents()
$ws.Range("S201").Value = 'chemotaxis (purine?), flagellum'
$ws.Range("R201").ClearContents()
$ws.Range("S203").Value = 'chemotaxis (purine?)'
$ws.Range("R203").ClearContents()

# --- Match the author's final on-screen selection ---
$ws.Range("M174").Select()
